$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs(1)

# Create a fresh empty paragraph right after the title. It inherits the
# title's Heading1 style, so reset it to Normal (the document's default
# body style, which stores no explicit pStyle) before filling it in.
$titlePara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs(2)
$newPara.Range.ParagraphFormat.Style = "Normal"

# Replace the (now empty) paragraph's full range -- including its paragraph
# mark -- with the target run structure in one shot, via a raw WordML
# fragment. This gives exact control of run boundaries/formatting, matching
# a leading empty run, a bold "Meta description" run, and a plain run with
# the rest of the sentence.
$metaXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r/>
<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
<w:r><w:t>: Read our honest review of Big Scary Fortune, a new Halloween slot game from Inspired Gaming. Play for free and experience its special features.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$newParaFullRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newParaFullRange.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Near the end of the document, drop the duplicated title paragraph and
#    turn the old meta-description sentence into the new image prompt.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
$dupTitleRange = $d.Range($dupTitlePara.Range.Start, $dupTitlePara.Range.End)
$dupTitleRange.Delete()

# Scope the find/replace to just the final paragraph's range so the earlier
# "Meta description: Read our honest review..." sentence (which shares the
# same wording) is left untouched.
$finalCount = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($finalCount)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End)

$finalRange.Find.Execute(
    "Read our honest review of Big Scary Fortune, a new Halloween slot game from Inspired Gaming. Play for free and experience its special features.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Prompt: Create a feature image for Big Scary Fortune that features a happy Maya warrior with glasses in cartoon style. The image should be eye-catching and show the warrior interacting with the spooky elements of the game, such as ghosts and pumpkins. The background should be dark and eerie to match the theme of the game. Make sure that the Maya warrior looks excited and ready to win big while playing Big Scary Fortune.",
    2
) | Out-Null
